# Scheduled-runner refresh of the Leve-profit market data across all
# job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Updates pull in
# newer currentAveragePrice(NQ/HQ) + LevePrice(NQ/HQ) + LeveProfit(NQ/HQ)
# numbers for the affected leves; a couple of rows gain/lose an HQ
# profit column entirely because HQ market data became available /
# unavailable for that item.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 788.84375
$ws.Range("I15").Value = 788.84375
$ws.Range("K15").Value = 2366.53125
$ws.Range("M15").Value = -2197.53125

$ws.Range("H48").Value = 7000
$ws.Range("J48").Value = 7000
$ws.Range("L48").Value = 21000
$ws.Range("N48").Value = -21584

$ws.Range("H56").Value = 7000
$ws.Range("J56").Value = 7000
$ws.Range("L56").Value = 21000
$ws.Range("N56").Value = -22068

$ws.Range("H62").Value = 1999
$ws.Range("I62").Value = 999
$ws.Range("K62").Value = 999
$ws.Range("M62").Value = -375

$ws.Range("H65").Value = 1999
$ws.Range("I65").Value = 999
$ws.Range("K65").Value = 4995
$ws.Range("M65").Value = -1875

$ws.Range("H92").Value = 791.0909
$ws.Range("I92").Value = 791.0909
$ws.Range("K92").Value = 791.0909
$ws.Range("M92").Value = 456.9091

$ws.Range("H138").Value = 5097.8823
$ws.Range("J138").Value = 9045.412
$ws.Range("L138").Value = 27136.236
$ws.Range("N138").Value = -37416.236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 63407.418
$ws.Range("I31").Value = 11266.625
$ws.Range("K31").Value = 11266.625
$ws.Range("M31").Value = -10972.625

$ws.Range("H32").Value = 5755.057
$ws.Range("I32").Value = 5336.0884
$ws.Range("K32").Value = 5336.0884
$ws.Range("M32").Value = -5049.0884

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H45").Value = 3239.4443
$ws.Range("I45").Value = 3025.75
$ws.Range("K45").Value = 3025.75
$ws.Range("M45").Value = -2648.75

$ws.Range("H74").Value = 26644.805
$ws.Range("I74").Value = 26644.805
$ws.Range("K74").Value = 26644.805
$ws.Range("M74").Value = -25770.805

$ws.Range("H77").Value = 26644.805
$ws.Range("I77").Value = 26644.805
$ws.Range("K77").Value = 133224.025
$ws.Range("M77").Value = -128856.025

$ws.Range("H102").Value = 4773.25
$ws.Range("I102").Value = 4169.4287
$ws.Range("J102").Value = 9000
$ws.Range("K102").Value = 4169.4287
$ws.Range("L102").Value = 9000
$ws.Range("M102").Value = -2547.4287
$ws.Range("N102").Value = -12244

$ws.Range("H122").Value = 16669666
$ws.Range("I122").Value = 2933.1667
$ws.Range("J122").Value = 33336400
$ws.Range("K122").Value = 8799.500100000001
$ws.Range("L122").Value = 100009200
$ws.Range("M122").Value = -6349.500100000001
$ws.Range("N122").Value = -100014100

$ws.Range("H132").Value = 107770.36
$ws.Range("I132").Value = 8823.75
$ws.Range("J132").Value = 701450
$ws.Range("K132").Value = 26471.25
$ws.Range("L132").Value = 2104350
$ws.Range("M132").Value = -23941.25
$ws.Range("N132").Value = -2109410

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 4226.6665
$ws.Range("I10").Value = 2562
$ws.Range("J10").Value = 12550
$ws.Range("K10").Value = 2562
$ws.Range("L10").Value = 12550
$ws.Range("M10").Value = -2422
$ws.Range("N10").Value = -12830

$ws.Range("H20").Value = 3235.9092
$ws.Range("I20").Value = 1867.6666
$ws.Range("K20").Value = 1867.6666
$ws.Range("M20").Value = -1620.6666

$ws.Range("H24").Value = 6005.3335
$ws.Range("I24").Value = 4008
$ws.Range("J24").Value = 10000
$ws.Range("K24").Value = 4008
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = -3773
$ws.Range("N24").Value = -10470

$ws.Range("H29").Value = 208
$ws.Range("I29").Value = 208
$ws.Range("K29").Value = 208
$ws.Range("M29").Value = 81

$ws.Range("H95").Value = 104912.25
$ws.Range("J95").Value = 104912.25
$ws.Range("L95").Value = 104912.25
$ws.Range("N95").Value = -110404.25

$ws.Range("H97").Value = 76408.8
$ws.Range("I97").Value = 22182.143
$ws.Range("K97").Value = 22182.143
$ws.Range("M97").Value = -21191.143

$ws.Range("H107").Value = 2420.7144
$ws.Range("I107").Value = 2146.7
$ws.Range("K107").Value = 2146.7
$ws.Range("M107").Value = -226.6999999999998

$ws.Range("H134").Value = 2785.3333
$ws.Range("I134").Value = 2942.4
$ws.Range("K134").Value = 8827.200000000001
$ws.Range("M134").Value = -6292.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2016.1666
$ws.Range("I16").Value = 2067.8462
$ws.Range("J16").Value = 1881.8
$ws.Range("K16").Value = 2067.8462
$ws.Range("L16").Value = 1881.8
$ws.Range("M16").Value = -1780.8462
$ws.Range("N16").Value = -2455.8

$ws.Range("H21").Value = 1444
$ws.Range("J21").Value = 1888
$ws.Range("L21").Value = 1888
$ws.Range("N21").Value = -2358

$ws.Range("H31").Value = 8624110
$ws.Range("I31").Value = 2646.6
$ws.Range("J31").Value = 27782918
$ws.Range("K31").Value = 2646.6
$ws.Range("L31").Value = 27782918
$ws.Range("M31").Value = -2351.6
$ws.Range("N31").Value = -27783508

$ws.Range("H34").Value = 8624110
$ws.Range("I34").Value = 2646.6
$ws.Range("J34").Value = 27782918
$ws.Range("K34").Value = 2646.6
$ws.Range("L34").Value = 27782918
$ws.Range("M34").Value = -2444.6
$ws.Range("N34").Value = -27783322

$ws.Range("H107").Value = 650.3077
$ws.Range("I107").Value = 630.7
$ws.Range("K107").Value = 630.7
$ws.Range("M107").Value = 1289.3

$ws.Range("H113").Value = 2016.1666
$ws.Range("I113").Value = 2067.8462
$ws.Range("J113").Value = 1881.8
$ws.Range("K113").Value = 2067.8462
$ws.Range("L113").Value = 1881.8
$ws.Range("M113").Value = 102.1538
$ws.Range("N113").Value = -6221.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 5476.273
$ws.Range("I9").Value = 33832.332
$ws.Range("K9").Value = 101496.996
$ws.Range("M9").Value = -101272.996

$ws.Range("H104").Value = 6990.8
$ws.Range("J104").Value = 6990.8
$ws.Range("L104").Value = 20972.4
$ws.Range("N104").Value = -26214.4

$ws.Range("H134").Value = 1081.875
$ws.Range("I134").Value = 815.1818
$ws.Range("K134").Value = 2445.5454
$ws.Range("M134").Value = 2624.4546

$ws.Range("H140").Value = 1576.3889
$ws.Range("I140").Value = 1269.3125
$ws.Range("K140").Value = 3807.9375
$ws.Range("M140").Value = 1372.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 12334.167
$ws.Range("I18").Value = 4005
$ws.Range("K18").Value = 4005
$ws.Range("M18").Value = -3712

$ws.Range("H70").Value = 6634.8335
$ws.Range("I70").Value = 4616.615
$ws.Range("K70").Value = 4616.615
$ws.Range("M70").Value = -4346.615

$ws.Range("H73").Value = 6634.8335
$ws.Range("I73").Value = 4616.615
$ws.Range("K73").Value = 4616.615
$ws.Range("M73").Value = -3680.615

$ws.Range("H80").Value = 5126.1665
$ws.Range("I80").Value = 3961.2856
$ws.Range("J80").Value = 5867.4546
$ws.Range("K80").Value = 3961.2856
$ws.Range("L80").Value = 5867.4546
$ws.Range("M80").Value = -2963.2856
$ws.Range("N80").Value = -7863.4546

$ws.Range("H83").Value = 5126.1665
$ws.Range("I83").Value = 3961.2856
$ws.Range("J83").Value = 5867.4546
$ws.Range("K83").Value = 19806.428
$ws.Range("L83").Value = 29337.273
$ws.Range("M83").Value = -14814.428
$ws.Range("N83").Value = -39321.273

$ws.Range("H97").Value = 1464.6875
$ws.Range("I97").Value = 764.1429000000001
$ws.Range("K97").Value = 764.1429000000001
$ws.Range("M97").Value = -268.1429000000001

$ws.Range("H132").Value = 2917.9092
$ws.Range("I132").Value = 2066.3333
$ws.Range("K132").Value = 6198.999899999999
$ws.Range("M132").Value = -3668.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4652.6665
$ws.Range("I7").Value = 3483
$ws.Range("K7").Value = 3483
$ws.Range("M7").Value = -3371

$ws.Range("H22").Value = 1642.1538
$ws.Range("I22").Value = 1196
$ws.Range("K22").Value = 1196
$ws.Range("M22").Value = -901

$ws.Range("H27").Value = 1642.1538
$ws.Range("I27").Value = 1196
$ws.Range("K27").Value = 1196
$ws.Range("M27").Value = -1089

$ws.Range("H55").Value = 650.42426
$ws.Range("I55").Value = 555
$ws.Range("K55").Value = 555
$ws.Range("M55").Value = -382

$ws.Range("H82").Value = 2695.1
$ws.Range("I82").Value = 2509.8
$ws.Range("K82").Value = 2509.8
$ws.Range("M82").Value = -2148.8

$ws.Range("H85").Value = 2695.1
$ws.Range("I85").Value = 2509.8
$ws.Range("K85").Value = 2509.8
$ws.Range("M85").Value = -1261.8

$ws.Range("H101").Value = 39393
$ws.Range("J101").Value = 39393
$ws.Range("L101").Value = 39393
$ws.Range("N101").Value = -45883

$ws.Range("H122").Value = 4314489
$ws.Range("I122").Value = 4085.5908
$ws.Range("J122").Value = 17861472
$ws.Range("K122").Value = 12256.7724
$ws.Range("L122").Value = 53584416
$ws.Range("M122").Value = -9806.7724
$ws.Range("N122").Value = -53589316

$ws.Range("H126").Value = 4652.6665
$ws.Range("I126").Value = 3483
$ws.Range("K126").Value = 10449
$ws.Range("M126").Value = -7979

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2314.2454
$ws.Range("I132").Value = 2174.578
$ws.Range("J132").Value = 3099.875
$ws.Range("K132").Value = 6523.734
$ws.Range("L132").Value = 9299.625
$ws.Range("M132").Value = -3993.734
$ws.Range("N132").Value = -14359.625
